$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering of model names for column A (rows 2..26)
$names = @(
    "model_11_1_0",
    "model_11_1_22",
    "model_11_1_21",
    "model_11_1_20",
    "model_11_1_19",
    "model_11_1_18",
    "model_11_1_17",
    "model_11_1_16",
    "model_11_1_15",
    "model_11_1_14",
    "model_11_1_13",
    "model_11_1_23",
    "model_11_1_12",
    "model_11_1_10",
    "model_11_1_9",
    "model_11_1_8",
    "model_11_1_7",
    "model_11_1_6",
    "model_11_1_5",
    "model_11_1_4",
    "model_11_1_3",
    "model_11_1_2",
    "model_11_1_1",
    "model_11_1_11",
    "model_11_1_24"
)

# Same metric values now shared by every data row
$values = @(0.3494677884409869, 0.4003320152976594, 0.2037968952667155, 0.3509831259012421, 0.7199474573135376, 0.9894625544548035, 0.8652127385139465, 0.930992066860199)

for ($i = 0; $i -lt $names.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $names[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $c = $j + 2
        $ws.Cells.Item($r, $c).Value = $values[$j]
    }
}
